$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.797.01'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").Value = '1.889.41'
$ws.Range("E3").Value = '  +1.15%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.94%  '

$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4800'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.89%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2957'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.81'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '101.31'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +18.70%  '

$ws.Range("D12").Value = '1.894.65'
$ws.Range("E12").Value = '  +1.43%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07622'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.146'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6596'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '306.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +26.70%  '

$ws.Range("D17").Value = '30.772.92'
$ws.Range("E17").Value = '  +0.57%  '

$ws.Range("E18").Value = '  +2.50%  '

$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007622'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.14%  '

$ws.Range("D21").Value = '2.134.22'
$ws.Range("E21").Value = '  +1.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.173'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.194'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.323'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.46%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.966'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.76%  '

$ws.Range("E29").Value = '  +9.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.347'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.195'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.008'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05085'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7437'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.165'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.54%  '

$ws.Range("E36").Value = '  +0.52%  '

$ws.Range("E37").Value = '  +3.76%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.710'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.057'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '109.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.34%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8820'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4213'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.657'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.84%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '68.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.08%  '

$ws.Range("E46").Value = '  -1.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.143'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.02%  '

$ws.Range("E48").Value = '  +0.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.95'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.59%  '

$ws.Range("E50").Value = '  +1.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.404'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.01%  '
